$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $ok = $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Host "MISS:" $old
    }
}

# 1) "负责增值物流、履约..." -> "负责物流计费、售后履约..."
Replace-Text "负责增值物流、履约、金融消费等增值项目，为公司提升毛利和收入" "负责物流计费、售后履约、金融消费等增值项目，为公司提升毛利和收入"

# 2) "为消费者提供退换货后的包运费投保、申请、理赔服务..." -> drop "、申请"
Replace-Text "为消费者提供退换货后的包运费投保、申请、理赔服务，并收取商家相应的服务费。" "为消费者提供退换货后的包运费投保、理赔服务，并收取商家相应的服务费。"

# 3) "搭建服务单、保单、理赔等领域服务..." -> "搭建服务、投保、理赔等领域能力..."
Replace-Text "搭建服务单、保单、理赔等领域服务，利用其上下文的解耦和内聚进一步解决理赔补偿、组件开放、保险抽象等难题" "搭建服务、投保、理赔等领域能力，利用核心领域的内聚和解耦进一步解决理赔补偿、组件开放、保险抽象等难题"

# 4) "...其核心都是通过服务赚取..." -> "...其核心都是提供服务赚取..."
Replace-Text "物流增值业务包括有赞寄件（面向商家）和上门取件（面向消费者），其核心都是通过服务赚取与三方物流间的差价" "物流增值业务包括有赞寄件（面向商家）和上门取件（面向消费者），其核心都是提供服务赚取与三方物流间的差价"

# 5) "有赞寄件是正向交易下单后物流发货、结算的核心，其通过物流商运营..." -> reworded
Replace-Text "有赞寄件是正向交易下单后物流发货、结算的核心，其通过物流商运营、运费定价等构造整个物流计费、结算模型" "有赞寄件是正向交易下单后物流发货、运费结算的核心，通过物流商运营、运费定价等构造整个计费、结算模型"

# 6) "上门取件是交易逆向售后的核心服务，其状态机 + " -> "...关键服务，其利用状态机 + "
Replace-Text "上门取件是交易逆向售后的核心服务，其状态机 + " "上门取件是交易逆向售后的关键服务，其利用状态机 + "

# 7) "架构驱动支付..." -> "策略驱动支付..."
Replace-Text "架构驱动支付、取件单、三方物流、交易单等状态一致" "策略驱动支付、取件单、三方物流、交易单等状态一致"

# 8) "模型上，通过业务组件可插拔..." paragraph rewritten (bookmark removed from here implicitly if present - it isn't here originally)
Replace-Text "模型上，通过业务组件可插拔 + 数据倒置依赖构建出抽象稳定的模型，实现对外开放能力，并完成业务的二次增长" "并在模型上，通过组件可插拔 + 数据倒置依赖构建出抽象稳定的模型，实现对外开放，完成业务的二次增长"

# 9) "基础设施上，搭建了与三方交互的物流体系，并通过心跳、监控、限流" + "等方式维系其稳定" (with _GoBack bookmark
#    straddling the two original runs) gets merged into a single run. Replacing the full text (which
#    spans across the embedded bookmark) removes the old _GoBack bookmark as a side effect.
Replace-Text "基础设施上，搭建了与三方交互的物流体系，并通过心跳、监控、限流等方式维系其稳定" "基础设施上，搭建了与三方交互的物流体系，并通过心跳、监控、限流等方式维系其稳定"

# 10) Re-create the _GoBack bookmark at its new location: right before "完成业务的二次增长"
#     in the now-rewritten "模型上" paragraph (i.e. right after "...实现对外开放，").
$find2 = $d.Content.Find
$find2.ClearFormatting()
$found2 = $find2.Execute("并在模型上，通过组件可插拔 + 数据倒置依赖构建出抽象稳定的模型，实现对外开放，", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $gb = $d.Range($find2.Parent.End, $find2.Parent.End)
    $d.Bookmarks.Add("_GoBack", $gb)
} else {
    Write-Host "Could not find anchor text for _GoBack bookmark"
}
